$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 658; everything currently at/after row 658
# (through the old last row 699) shifts down to rows 660..701.
$ws.Rows.Item(658).Resize(2).Insert()

# The two newly-inserted rows hold brand-new data points for 2026/01/17
# (Sat) and 2026/01/18 (Sun). Force text formatting first so the
# date-like strings and weekday kanji are stored as literal text (as in
# the rest of the column) instead of being auto-converted to date
# serials.
$ws.Range("A658:B659").NumberFormat = "@"

$ws.Range("A658").Value = "2026/01/17"
$ws.Range("B658").Value = "土"
$ws.Range("C658").Value = 22
$ws.Range("D658").Value = 201

$ws.Range("A659").Value = "2026/01/18"
$ws.Range("B659").Value = "日"
$ws.Range("C659").Value = 2
$ws.Range("D659").Value = 201

# Restore the default (General/Normal) style so the new cells match the
# unstyled look of every other data row instead of keeping the
# temporary "Text" number format applied above.
$ws.Range("A658:D659").Style = "Normal"
